$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ひとこと" (B) column text for each entry was re-typed/pasted so that each
# piece now reads as a single run-on paragraph instead of several lines broken
# up by newlines. Cells are touched in the same order the new shared strings
# appear in the saved workbook (B8, B7, B9, B6, B5, B2, B3, B4).
$ws.Cells.Item(8,2).Value = "有給休暇を取ったある日、わたしは車を走らせていた。`n雨が降るなか､立ち寄ったケンタッキーフライドチキンで､" + [char]0x0B + "チキンとポテトを買い、車に戻ってラジオをかけて食べた。誰にも気を遣わず､何の義務もなく､ただ雨とラジオの音とチキンの香りに包まれた時間。"
$ws.Cells.Item(7,2).Value = "休日の前夜、燻した６Pチーズと堅あげポテト" + [char]0x3000 + "下戸なのに、ウイスキーロック" + [char]0x3000 + "いつのまにか、夢の中"
$ws.Cells.Item(9,2).Value = "何も買わない。何も食べない。何も生まれない。ただ、流れるその時を伸ばして、歩きつづける。休日の昼を、所有している。LINEはときどき鳴っている。妻や子供のことを思い出しながら、新書に手を伸ばす。"
$ws.Cells.Item(6,2).Value = "レモンサワーを頼む。普段はレモンは１個だけしか絞らないけど、今日は2個絞ろうかなぁ。だって今日は､うるめが美味しそうだから。"
$ws.Cells.Item(5,2).Value = "毎週日曜日、お酒を飲みながら、ひとりで「ちびまる子ちゃん」を観る。"
$ws.Cells.Item(2,2).Value = "特に予定のない土日、家の近くの公園でランニングをしたり、サウナに行ってととのう時間。よりエネルギッシュになるために、あえてリラックス・リフレッシュする意識を持つ。"
$ws.Cells.Item(3,2).Value = "週末の深夜、家の前を流れる川の堤防で一服。誰もない、誰からも干渉されない自分だけの時間。この一本が、たまらなくうまい。"
$ws.Cells.Item(4,2).Value = "知らない土地に出張したとき、誰も自分を知らない町を歩く。見知らぬ人々の営みを思い描きながら歩くその瞬間に、不思議な楽しさを感じる。"

# Leave the selection on B7, matching the saved cursor position.
$ws.Range("B7").Select()
